$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source values are numeric-looking strings that must remain stored as
# text (as in the original workbook). Prefixing with a leading apostrophe
# forces Excel to keep them as text instead of auto-converting to numbers.
$ws.Range("B10").Value = "'20.11"
$ws.Range("C10").Value = "'1.37"
$ws.Range("D10").Value = "'21.49"

$ws.Range("C12").Value = "'21.76"
$ws.Range("D12").Value = "'84.36"

$ws.Range("B14").Value = "'93.46"
$ws.Range("C14").Value = "'6.39"
$ws.Range("D14").Value = "'99.84"

$ws.Range("D29").Value = "'19.37"
